$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '35.583.68'
$ws.Range('E2').Value = '  +1.46%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.916.62'
$ws.Range('E3').Value = '  +3.41%  '

$ws.Range('E4').Value = '  +0.52%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '247.24'
$ws.Range('E5').Value = '  +4.26%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.654'
$ws.Range('E6').Value = '  +5.16%  '

$ws.Range('E7').Value = '  +0.45%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '42.22'
$ws.Range('E8').Value = '  -0.17%  '

$ws.Range('E9').Value = '  +6.03%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '49.26'
$ws.Range('E10').Value = '  +5.43%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0727'
$ws.Range('E11').Value = '  +4.68%  '

$ws.Range('E12').Value = '  +1.08%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.196.65'

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '12.36'
$ws.Range('E14').Value = '  +8.02%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.702'
$ws.Range('E15').Value = '  +3.53%  '

$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.95'
$ws.Range('E16').Value = '  +5.21%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.919.23'
$ws.Range('E17').Value = '  +3.78%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '35.611.46'
$ws.Range('E18').Value = '  +1.65%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '72.50'
$ws.Range('E19').Value = '  +3.15%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0829'
$ws.Range('E20').Value = '  +4.15%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '246.25'
$ws.Range('E21').Value = '  +2.38%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '12.70'
$ws.Range('E22').Value = '  +3.79%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.88'
$ws.Range('E23').Value = '  +2.02%  '

$ws.Range('E24').Value = '  +0.49%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.32'
$ws.Range('E25').Value = '  +2.10%  '

$ws.Range('E26').Value = '  +16.43%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '171.36'
$ws.Range('E27').Value = '  +0.35%  '

$ws.Range('E28').Value = '  +6.61%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '18.54'
$ws.Range('E29').Value = '  +4.86%  '

$ws.Range('E30').Value = '  +3.01%  '

$ws.Range('E31').Value = '  +4.28%  '

$ws.Range('E32').Value = '  +2.79%  '

$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.21'
$ws.Range('E33').Value = '  +3.99%  '

$ws.Range('E34').Value = '  +0.55%  '

$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.937'
$ws.Range('E35').Value = '  +19.46%  '

$ws.Range('E36').Value = '  +5.68%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.04'
$ws.Range('E37').Value = '  +1.86%  '

$ws.Range('E38').Value = '  +1.49%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.11'
$ws.Range('E39').Value = '  +2.72%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0212'
$ws.Range('E40').Value = '  +5.04%  '

$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0642'
$ws.Range('E41').Value = '  +16.36%  '

$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '91.95'
$ws.Range('E42').Value = '  +0.93%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '15.84'
$ws.Range('E43').Value = '  +7.44%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.359.42'
$ws.Range('E44').Value = '  +0.37%  '

$ws.Range('E45').Value = '  +2.35%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '46.98'
$ws.Range('E46').Value = '  +36.48%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '12.66'
$ws.Range('E47').Value = '  +0.42%  '

$ws.Range('E48').Value = '  +2.35%  '

$ws.Range('E49').Value = '  -0.21%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.58'
$ws.Range('E50').Value = '  +0.84%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.106.97'
$ws.Range('E51').Value = '  +3.53%  '

Write-Host "Applied crypto price updates"